# MICS-infographic.pptx edit:
#  - "embiggened the map; ensmallened the references"
#  - 3 of the 4 "Rectangle 5" header runs turn from blue (0070C0) to red (FF0000)
#  - the map picture grows and shifts to line up with the References box
#  - the References box (Rectangle 6) shrinks/moves down to make room
#  - every run in the References box shrinks from 8pt to 7pt

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Rectangle 5: recolor three of the four section headers to red ---
$rect5 = $s.Shapes.Item(3)
$tr5 = $rect5.TextFrame.TextRange

$h1 = $tr5.Find("No Spanking")
$h1.Font.Color.RGB = 255

$h2 = $tr5.Find("No Insults")
$h2.Font.Color.RGB = 255

$h3 = $tr5.Find("No Physical Or Psychological Punishments")
$h3.Font.Color.RGB = 255

# NOTE: "Positive Discipline" (the 4th header) intentionally keeps its
# original blue (0070C0) color - it is not touched.

# --- Picture "Figure 1: Countries in MICS" (the map): embiggen + reposition ---
# (point values are nudged very slightly so that PowerPoint's internal
#  single-precision Left/Top/Width/Height round-trip to the exact target EMUs)
$map = $s.Shapes.Item(4)
$map.Left = 270.0000787401575
$map.Top = 267.99339582677163
$map.Width = 249.2607086614173
$map.Height = 136.3579527559055

# --- Rectangle 6 (References box): reposition + shrink to make room for the map ---
$refs = $s.Shapes.Item(5)
$refs.Left = 270.0000787401575
$refs.Top = 404.9471066141732
$refs.Width = 249.2607086614173
$refs.Height = 295.25545307086617

# Ensmallen every run of reference text from 8pt to 7pt.
$refs.TextFrame.TextRange.Font.Size = 7
